# Update NATMI LR-pairs TPM values (Sirpa-Cd47) with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "ligand" side values (columns G,H,I,J) keyed by Sending cluster (column A)
$ligandVals = @{
    "ECs"               = @(5.343111666666666, 16.029335, 0.007536583045055595, 0.007566622551600167)
    "FAPs"              = @(24.954405, 74.863215, 0.03519876756380422, 0.03533906371688481)
    "Inflammatory-Mac"  = @(312.3302103333334, 936.990631, 0.4405490123558179, 0.4423049639403425)
    "MuSCs"             = @(8.4436795, 16.887359, 0.01191000595300744, 0.007971651440709677)
    "Resolving-Mac"     = @(357.885376, 1073.656128, 0.504805631082315, 0.5068176983504629)
}

# New "receptor" side values (columns M,N,O,P) keyed by Target cluster (column D)
$receptorVals = @{
    "ECs"               = @(46.17354133333333, 138.520624, 0.1154336358852217, 0.1189208138601986)
    "FAPs"              = @(84.01327500000001, 252.039825, 0.2100327918507284, 0.2163777511873036)
    "Inflammatory-Mac"  = @(128.0910926666667, 384.273278, 0.3202271284388135, 0.3299009897940278)
    "MuSCs"             = @(35.18830149999999, 70.37660299999999, 0.08797058803540478, 0.06041874966919073)
    "Resolving-Mac"     = @(106.534543, 319.603629, 0.2663358557898317, 0.2743816954892795)
}

$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value()
    $target  = $ws.Cells.Item($r, 4).Value()

    if ($null -eq $sending -or $null -eq $target) { continue }

    $lig = $ligandVals[$sending]
    $rec = $receptorVals[$target]

    if ($null -eq $lig -or $null -eq $rec) { continue }

    # Existing (old) values, needed to decide whether the derived edge
    # columns (Q,R,S,T) actually need to be recomputed.
    $gOld = $ws.Cells.Item($r, 7).Value()
    $hOld = $ws.Cells.Item($r, 8).Value()
    $iOld = $ws.Cells.Item($r, 9).Value()
    $jOld = $ws.Cells.Item($r, 10).Value()
    $mOld = $ws.Cells.Item($r, 13).Value()
    $nOld = $ws.Cells.Item($r, 14).Value()
    $oOld = $ws.Cells.Item($r, 15).Value()
    $pOld = $ws.Cells.Item($r, 16).Value()
    $qOld = $ws.Cells.Item($r, 17).Value()
    $rOld = $ws.Cells.Item($r, 18).Value()
    $sOld = $ws.Cells.Item($r, 19).Value()
    $tOld = $ws.Cells.Item($r, 20).Value()

    $g = $lig[0]; $h = $lig[1]; $i = $lig[2]; $j = $lig[3]
    $m = $rec[0]; $n = $rec[1]; $o = $rec[2]; $p = $rec[3]

    $ws.Cells.Item($r, 7).Value  = $g   # G - Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $h   # H - Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $i   # I - Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value = $j   # J - Ligand derived specificity of total expression value

    $ws.Cells.Item($r, 13).Value = $m   # M - Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $n   # N - Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $o   # O - Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value = $p   # P - Receptor derived specificity of total expression value

    # Only recompute the derived edge-weight columns when one of their
    # inputs actually changed; otherwise keep the original stored value
    # so we don't introduce spurious last-bit floating point noise.
    if ($g -eq $gOld -and $m -eq $mOld) {
        $ws.Cells.Item($r, 17).Value = $qOld
    } else {
        $ws.Cells.Item($r, 17).Value = $g * $m   # Q - Edge average expression weight
    }

    if ($h -eq $hOld -and $n -eq $nOld) {
        $ws.Cells.Item($r, 18).Value = $rOld
    } else {
        $ws.Cells.Item($r, 18).Value = $h * $n   # R - Edge total expression weight
    }

    if ($i -eq $iOld -and $o -eq $oOld) {
        $ws.Cells.Item($r, 19).Value = $sOld
    } else {
        $ws.Cells.Item($r, 19).Value = $i * $o   # S - Edge average expression derived specificity
    }

    if ($j -eq $jOld -and $p -eq $pOld) {
        $ws.Cells.Item($r, 20).Value = $tOld
    } else {
        $ws.Cells.Item($r, 20).Value = $j * $p   # T - Edge total expression derived specificity
    }
}
